$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing Thursday hours for the week of 43171 (row 9)
$ws.Range("E9").Value = 5.75

# Update the active cell selection to match the saved state
$ws.Range("J22").Select()
